$wb = $excel.ActiveWorkbook

# ---- Section_A ----
$wsA = $wb.Worksheets.Item("Section_A")
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "EC264"
$wsA.Range("B3").Value = "Free"
$wsA.Range("C3").Value = "Free"
$wsA.Range("E3").Value = "EC304"
$wsA.Range("B5").Value = "Free"
$wsA.Range("C5").Value = "CS307"
$wsA.Range("D5").Value = "EC303"
$wsA.Range("E5").Value = "CS307"
$wsA.Range("F5").Value = "CS307"
$wsA.Range("D6").Value = "EC304 (Tutorial)"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "Free"
$wsA.Range("B7").Value = "EC264"
$wsA.Range("C7").Value = "EC303"
$wsA.Range("D7").Value = "EC262"
$wsA.Range("E7").Value = "Free"
$wsA.Range("F7").Value = "EC264"
$wsA.Range("D8").Value = "CS307 (Tutorial)"

# ---- Section_B ----
$wsB = $wb.Worksheets.Item("Section_B")
$wsB.Range("B2").Value = "EC262"
$wsB.Range("C2").Value = "EC304"
$wsB.Range("D2").Value = "CS307"
$wsB.Range("E2").Value = "EC303"
$wsB.Range("F2").Value = "EC264"
$wsB.Range("B3").Value = "EC264"
$wsB.Range("C3").Value = "Free"
$wsB.Range("D3").Value = "EC303"
$wsB.Range("E3").Value = "Free"
$wsB.Range("F3").Value = "CS307"
$wsB.Range("C5").Value = "Free"
$wsB.Range("E5").Value = "EC304"
$wsB.Range("F5").Value = "EC303"
$wsB.Range("D6").Value = "EC304 (Tutorial)"
$wsB.Range("B7").Value = "Free"
$wsB.Range("E7").Value = "EC262"
$wsB.Range("F7").Value = "EC262"
$wsB.Range("D8").Value = "Free"
